# header shrinking + docs + career loop
#
# Adds a new "Header Tag" column (E) populated with "Information Technology"
# for the existing program rows, and appends a new trailer/test row (row 7)
# with "Test data" / "test data" / "/test" used for a documentation /
# career-loop smoke-test entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "Header Tag" column (E) ------------------------------------------
$ws.Range("E1").Value = "Header Tag"
$ws.Range("E3").Value = "Information Technology"
$ws.Range("E4").Value = "Information Technology"
$ws.Range("E5").Value = "Information Technology"
$ws.Range("E6").Value = "Information Technology"

# --- New trailer row (7) ----------------------------------------------------
$ws.Range("A7").Value = "Test data"
$ws.Range("B7").Value = "test data"
$ws.Range("D7").Value = "/test"

# Match the wrap-text / vertical-centered formatting already used by column
# B and column D on the rows above by copying their formats down.
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View / selection state --------------------------------------------------
# Select the whole of the new column and scroll the window down to the new row.
$ws.Columns.Item(5).Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
